$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.100.86"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.526.88"
$ws.Range("E3").Value = "  -2.29%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.87"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "172.79"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.528"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "2.525.08"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("E12").Value = "  -0.41%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.343"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.70%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.58"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "2.984.64"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "66.937.05"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "2.528.08"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  +4.63%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.36"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.71%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "353.81"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.44%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.18"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +6.02%  "
$ws.Range("E25").Value = "  +0.03%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "69.61"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.95"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D30").Value = "0.0₃0977"
$ws.Range("E30").Value = "  -1.04%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "532.37"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.15"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("E33").Value = "  +0.09%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -0.82%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -0.08%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "157.81"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.62"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.354"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  -0.06%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.50"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.25%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "149.06"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "0.0₆0278"
$ws.Range("E48").Value = "  -2.90%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -1.04%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0757"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
